# Refactor item spawning logic and add default values to Item properties
#
# 1. Capture the current per-row Coin/Banana spawn amounts (B2:C17).
# 2. Overwrite B2:B17 / C2:C17 with the new default spawn values (0 / 200)
#    used by the refactored spawning logic.
# 3. Re-home the old per-row values as a static reference table in the new
#    rows 56:71 (columns C/D) so the original tuning numbers aren't lost.
# 4. Tidy up the sheet view (drop the pinned topLeftCell, move the
#    selection to R10) and narrow column O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17

# --- Preserve the old B2:C17 values before overwriting them ---
$oldB = @()
$oldC = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldB += , $ws.Cells.Item($r, 2).Value2
    $oldC += , $ws.Cells.Item($r, 3).Value2
}

# --- Set new default B/C values for rows 2-17 ---
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value2 = 0
    $ws.Cells.Item($r, 3).Value2 = 200
}

# --- Write old values into the new reference table at rows 56-71, columns C/D ---
$destRow = 56
for ($i = 0; $i -lt $oldB.Count; $i++) {
    $ws.Cells.Item($destRow, 3).Value2 = $oldB[$i]
    $ws.Cells.Item($destRow, 4).Value2 = $oldC[$i]
    $destRow++
}

# --- Narrow column O ---
$ws.Columns.Item(15).ColumnWidth = 18.5

# --- Update sheet view: drop the pinned top-left cell, move selection to R10 ---
$ws.Range("R10").Select()

# --- Reposition the workbook window ---
$excel.ActiveWindow.Left = -120
